$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}

Replace-Text "2024-12-29 Sunday" "2024-12-30 Monday"

Replace-Text "692×2=1384" "393×3=1179"
Replace-Text "305×2=610" "693×2=1386"
Replace-Text "128×6=768" "314×3=942"
Replace-Text "762×6=4572" "607×8=4856"
Replace-Text "374×8=2992" "505×8=4040"

Replace-Text "856×2=1712" "202×2=404"
Replace-Text "502×2=1004" "503×7=3521"
Replace-Text "753×3=2259" "816×4=3264"
Replace-Text "237×7=1659" "722×3=2166"
Replace-Text "414×4=1656" "947×6=5682"

Replace-Text "722×6=4332" "925×6=5550"
Replace-Text "705×8=5640" "160×4=640"
Replace-Text "473×6=2838" "342×2=684"
Replace-Text "702×6=4212" "681×5=3405"
Replace-Text "308×5=1540" "780×4=3120"

Replace-Text "754×8=6032" "962×2=1924"
Replace-Text "888×2=1776" "944×2=1888"
Replace-Text "739×5=3695" "307×6=1842"
Replace-Text "772×2=1544" "583×6=3498"
Replace-Text "863×3=2589" "109×9=981"

Replace-Text "388×2=776" "743×7=5201"
Replace-Text "613×8=4904" "118×7=826"
Replace-Text "616×8=4928" "797×5=3985"
Replace-Text "342×6=2052" "322×2=644"
Replace-Text "456×3=1368" "306×9=2754"
